# Corrected game details in final presentation
#
# Slide 2 ("The Project") — three small wording fixes inside the bullet body.
# Slide 3 ("Game Concept") — five bullet rewrites, also collapsing the blank
#   spacer paragraphs that used to sit between each bullet.

$p = $ppt.ActivePresentation

# ---- Slide 2: body bullet text fixes -------------------------------------
$slide2 = $p.Slides.Item(2)
$body2 = $slide2.Shapes.Item(2).TextFrame.TextRange

$body2.Paragraphs(1, 1).Text = "Puzzle Dots is a web-based game in which colored dots are manipulated to solve puzzles. "
$body2.Paragraphs(2, 1).Text = "The game offers short levels with varying complexity and replayability through finding improved solutions that require less moves. "
$body2.Paragraphs(3, 1).Text = "Written in Javascript and HTML5 with CSS."

# ---- Slide 3: bullet rewrites + removal of blank spacer paragraphs -------
$slide3 = $p.Slides.Item(3)
$body3 = $slide3.Shapes.Item(2).TextFrame.TextRange

$newBullets = @(
    "Pieces move simultaneously in groups, where each dot can have a different direction.",
    "When pieces move into each other, different reactions occur.",
    "Forming a new piece through a reaction is called blending.",
    "Reactions between pieces can be learned through experimentation or by reading the rules.",
    "The goal is to pair all pieces and spaces marked with goals by matching colors."
)

$body3.Text = [string]::Join("`r", $newBullets)
